# Initial Contour USB check-in. Header records reading and parsed.
#
# - Adds two new Bayer HealthCare meter names to column B of the
#   "Meters" sheet: "Contour Next USB*" (B7) and "Contour Next EZ*" (B8).
# - Applies the "Input" cell style to B5 (matching the styling already
#   used for similar cells elsewhere on the sheet).
# - Leaves the sheet's active selection on B13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meters")

$ws.Range("B5").Style = "Input"

$ws.Range("B7").Value = "Contour Next USB*"
$ws.Range("B8").Value = "Contour Next EZ*"

$ws.Range("B13").Select()
